$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @{ "A2"="ECs"; "B2"="a"; "C2"="F11r"; "D2"="ECs"; "E2"=1; "F2"=0.3333333333333333; "G2"=0.045339; "H2"=0.136017; "I2"=0.1740293637846656; "J2"=0.1740293637846656; "K2"=3; "L2"=1; "M2"=40.91514966666667; "N2"=122.745449; "O2"=0.8529192913871414; "P2"=0.8529192913871415; "Q2"=1.855051970737; "R2"=16.695467736633; "S2"=0.148433001639772; "T2"=0.148433001639772 },
    @{ "A3"="ECs"; "B3"="a"; "C3"="F11r"; "D3"="FAPs"; "E3"=1; "F3"=0.3333333333333333; "G3"=0.045339; "H3"=0.136017; "I3"=0.1740293637846656; "J3"=0.1740293637846656; "K3"=2; "L3"=0.6666666666666666; "M3"=0.165216; "N3"=0.495648; "O3"=0.00344410114086962; "P3"=0.003444101140869621; "Q3"=0.007490728223999999; "R3"=0.06741655401600001; "S3"=0.0005993747303555809; "T3"=0.0005993747303555812 },
    @{ "A4"="ECs"; "B4"="a"; "C4"="F11r"; "D4"="MuSCs"; "E4"=1; "F4"=0.3333333333333333; "G4"=0.045339; "H4"=0.136017; "I4"=0.1740293637846656; "J4"=0.1740293637846656; "K4"=2; "L4"=0.6666666666666666; "M4"=0.4441646666666667; "N4"=1.332494; "O4"=0.009259079236881667; "P4"=0.009259079236881667; "Q4"=0.020137981822; "R4"=0.181241836398; "S4"=0.001611351668826323; "T4"=0.001611351668826324 },
    @{ "A5"="ECs"; "B5"="a"; "C5"="F11r"; "D5"="Resolving-Mac"; "E5"=1; "F5"=0.3333333333333333; "G5"=0.045339; "H5"=0.136017; "I5"=0.1740293637846656; "J5"=0.1740293637846656; "K5"=3; "L5"=1; "M5"=6.446186333333333; "N5"=19.338559; "O5"=0.1343775282351073; "P5"=0.1343775282351073; "Q5"=0.292263642167; "R5"=2.630372779503; "S5"=0.02338563574571166; "T5"=0.02338563574571166 },
    @{ "A6"="MuSCs"; "B6"="a"; "C6"="F11r"; "D6"="ECs"; "E6"=1; "F6"=0.3333333333333333; "G6"=0.215186; "H6"=0.645558; "I6"=0.8259706362153344; "J6"=0.8259706362153345; "K6"=3; "L6"=1; "M6"=40.91514966666667; "N6"=122.745449; "O6"=0.8529192913871414; "P6"=0.8529192913871415; "Q6"=8.804367396171333; "R6"=79.239306565542; "S6"=0.7044862897473694; "T6"=0.7044862897473696 },
    @{ "A7"="MuSCs"; "B7"="a"; "C7"="F11r"; "D7"="FAPs"; "E7"=1; "F7"=0.3333333333333333; "G7"=0.215186; "H7"=0.645558; "I7"=0.8259706362153344; "J7"=0.8259706362153345; "K7"=2; "L7"=0.6666666666666666; "M7"=0.165216; "N7"=0.495648; "O7"=0.00344410114086962; "P7"=0.003444101140869621; "Q7"=0.035552170176; "R7"=0.319969531584; "S7"=0.00284472641051404; "T7"=0.00284472641051404 },
    @{ "A8"="MuSCs"; "B8"="a"; "C8"="F11r"; "D8"="MuSCs"; "E8"=1; "F8"=0.3333333333333333; "G8"=0.215186; "H8"=0.645558; "I8"=0.8259706362153344; "J8"=0.8259706362153345; "K8"=2; "L8"=0.6666666666666666; "M8"=0.4441646666666667; "N8"=1.332494; "O8"=0.009259079236881667; "P8"=0.009259079236881667; "Q8"=0.09557801796133333; "R8"=0.860202161652; "S8"=0.007647727568055343; "T8"=0.007647727568055344 },
    @{ "A9"="MuSCs"; "B9"="a"; "C9"="F11r"; "D9"="Resolving-Mac"; "E9"=1; "F9"=0.3333333333333333; "G9"=0.215186; "H9"=0.645558; "I9"=0.8259706362153344; "J9"=0.8259706362153345; "K9"=3; "L9"=1; "M9"=6.446186333333333; "N9"=19.338559; "O9"=0.1343775282351073; "P9"=0.1343775282351073; "Q9"=1.387129052324667; "R9"=12.484161470922; "S9"=0.1109918924893956; "T9"=0.1109918924893957 }
)

foreach ($row in $rowsData) {
    foreach ($key in $row.Keys) {
        $ws.Range($key).Value = $row[$key]
    }
}